# Working on comp 240 Flask
# Add three more invoice sheets ("Invoice 2", "Invoice 3", "Invoice 4"),
# each a copy of "Invoice 1" with its own data, and move the selection
# on "Invoice 1" to H27.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Invoice 1")

function Fill-Invoice {
    param(
        $ws,
        [int]$num,
        [string]$customer_id,
        [string]$f_name,
        [string]$l_name,
        [string]$phone,
        [string]$invoice_date,
        [string]$payment_method,
        [string]$total_price
    )

    $ws.Range("A1").Value = "INVOICE #$num"

    # customer_id and total_price look numeric but the source data stores
    # them as plain text (like every other value column here) - force text
    # so they don't get auto-coerced into numbers.
    $ws.Range("C2").NumberFormat = "@"
    $ws.Range("C2").Value = $customer_id

    $ws.Range("C3").Value = "`"$f_name`""
    $ws.Range("C4").Value = "`"$l_name`""
    $ws.Range("C5").Value = "`"$phone`""
    $ws.Range("C6").Value = "`"$invoice_date`""
    $ws.Range("C7").Value = "`"$payment_method`""

    $ws.Range("C8").NumberFormat = "@"
    $ws.Range("C8").Value = $total_price
}

# Invoice 2 - copied right after Invoice 1
$ws1.Copy([System.Reflection.Missing]::Value, $ws1)
$ws2 = $wb.Worksheets.Item($ws1.Index + 1)
$ws2.Name = "Invoice 2"
Fill-Invoice $ws2 2 "2" "Richmound" "Hulmes" "941-402-4909" "2020-03-07" "Credit" "1282.5"

# Invoice 3 - copied right after Invoice 2
$ws2.Copy([System.Reflection.Missing]::Value, $ws2)
$ws3 = $wb.Worksheets.Item($ws2.Index + 1)
$ws3.Name = "Invoice 3"
Fill-Invoice $ws3 3 "3" "Annetta" "Colleford" "198-751-6091" "2020-03-08" "Debit" "210"

# Invoice 4 - copied right after Invoice 3
$ws3.Copy([System.Reflection.Missing]::Value, $ws3)
$ws4 = $wb.Worksheets.Item($ws3.Index + 1)
$ws4.Name = "Invoice 4"
Fill-Invoice $ws4 4 "4" "Jessalyn" "Formie" "103-907-2970" "2020-03-24" "Cheque" "464.5"

# Move the active selection on "Sheet1" to H27
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Range("H27").Select()
